$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '43.397.24'
Set-TextValue $ws.Range("E2") '  -6.88%  '
Set-TextValue $ws.Range("D3") '2.583.23'
Set-TextValue $ws.Range("E3") '  -0.96%  '
Set-TextValue $ws.Range("E4") '  +0.12%  '
Set-TextValue $ws.Range("D5") '299.75'
Set-TextValue $ws.Range("E5") '  -2.98%  '
Set-TextValue $ws.Range("D6") '95.79'
Set-TextValue $ws.Range("E6") '  -4.98%  '
Set-TextValue $ws.Range("D7") '0.575'
Set-TextValue $ws.Range("E7") '  -4.98%  '
Set-TextValue $ws.Range("E8") '  +0.19%  '
Set-TextValue $ws.Range("D9") '0.555'
Set-TextValue $ws.Range("E9") '  -4.67%  '
Set-TextValue $ws.Range("D10") '36.66'
Set-TextValue $ws.Range("E10") '  -7.32%  '
Set-TextValue $ws.Range("D11") '0.0812'
Set-TextValue $ws.Range("E11") '  -4.34%  '
Set-TextValue $ws.Range("D12") '7.79'
Set-TextValue $ws.Range("E12") '  -4.91%  '
Set-TextValue $ws.Range("D13") '2.988.06'
Set-TextValue $ws.Range("E13") '  -0.22%  '
Set-TextValue $ws.Range("E14") '  +0.89%  '
Set-TextValue $ws.Range("D15") '2.601.13'
Set-TextValue $ws.Range("E15") '  -0.04%  '
Set-TextValue $ws.Range("D16") '0.887'
Set-TextValue $ws.Range("E16") '  -4.47%  '
Set-TextValue $ws.Range("D17") '14.27'
Set-TextValue $ws.Range("E17") '  -5.47%  '
Set-TextValue $ws.Range("D18") '43.445.35'
Set-TextValue $ws.Range("E18") '  -6.88%  '
Set-TextValue $ws.Range("D19") '6.63'
Set-TextValue $ws.Range("E19") '  -2.67%  '
Set-TextValue $ws.Range("D20") '0.0₃0974'
Set-TextValue $ws.Range("E20") '  -4.73%  '
Set-TextValue $ws.Range("D21") '12.25'
Set-TextValue $ws.Range("E21") '  -6.73%  '
Set-TextValue $ws.Range("D22") '72.77'
Set-TextValue $ws.Range("E22") '  +1.18%  '
Set-TextValue $ws.Range("D23") '265.47'
Set-TextValue $ws.Range("E23") '  -4.44%  '
Set-TextValue $ws.Range("D24") '2.21'
Set-TextValue $ws.Range("E24") '  +0.83%  '
Set-TextValue $ws.Range("D25") '2.91'
Set-TextValue $ws.Range("E25") '  -4.88%  '
Set-TextValue $ws.Range("D26") '29.23'
Set-TextValue $ws.Range("E26") '  -0.30%  '
Set-TextValue $ws.Range("E27") '  -0.09%  '
Set-TextValue $ws.Range("D28") '10.22'
Set-TextValue $ws.Range("E28") '  -4.56%  '
Set-TextValue $ws.Range("D29") '2.21'
Set-TextValue $ws.Range("E29") '  -2.44%  '
Set-TextValue $ws.Range("D30") '37.28'
Set-TextValue $ws.Range("E30") '  -5.21%  '
Set-TextValue $ws.Range("D31") '6.01'
Set-TextValue $ws.Range("E31") '  -5.63%  '
Set-TextValue $ws.Range("D32") '3.57'
Set-TextValue $ws.Range("E32") '  -1.32%  '
Set-TextValue $ws.Range("D33") '2.23'
Set-TextValue $ws.Range("E33") '  +0.08%  '
Set-TextValue $ws.Range("D34") '151.82'
Set-TextValue $ws.Range("E34") '  -0.23%  '
Set-TextValue $ws.Range("E35") '  -2.08%  '
Set-TextValue $ws.Range("D36") '0.0811'
Set-TextValue $ws.Range("E36") '  -3.94%  '
Set-TextValue $ws.Range("E37") '  -6.00%  '
Set-TextValue $ws.Range("D38") '24.44'
Set-TextValue $ws.Range("E38") '  +4.78%  '
Set-TextValue $ws.Range("D39") '0.120'
Set-TextValue $ws.Range("E39") '  -2.13%  '
Set-TextValue $ws.Range("D40") '16.68'
Set-TextValue $ws.Range("E40") '  +2.42%  '
Set-TextValue $ws.Range("D41") '3.51'
Set-TextValue $ws.Range("E41") '  -4.09%  '
Set-TextValue $ws.Range("D42") '0.0314'
Set-TextValue $ws.Range("E42") '  -6.16%  '
Set-TextValue $ws.Range("D43") '3.83'
Set-TextValue $ws.Range("E43") '  -6.67%  '
Set-TextValue $ws.Range("D44") '2.058.77'
Set-TextValue $ws.Range("E44") '  -4.07%  '
Set-TextValue $ws.Range("E45") '  +0.14%  '
Set-TextValue $ws.Range("D46") '88.11'
Set-TextValue $ws.Range("E46") '  -5.65%  '
Set-TextValue $ws.Range("D47") '9.03'
Set-TextValue $ws.Range("E47") '  -5.50%  '
Set-TextValue $ws.Range("D48") '2.845.04'
Set-TextValue $ws.Range("E48") '  -0.26%  '
Set-TextValue $ws.Range("D49") '1.60'
Set-TextValue $ws.Range("E49") '  -0.04%  '
Set-TextValue $ws.Range("D50") '105.28'
Set-TextValue $ws.Range("E50") '  -4.28%  '
Set-TextValue $ws.Range("D51") '0.189'
Set-TextValue $ws.Range("E51") '  -6.28%  '
